$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.479.45'
$ws.Range("E2").Value = '  +1.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.689.42'
$ws.Range("E3").Value = '  +5.20%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.20'
$ws.Range("E5").Value = '  +4.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.16'
$ws.Range("E6").Value = '  +1.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.682.15'
$ws.Range("E9").Value = '  +4.84%  '

$ws.Range("E10").Value = '  +1.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.69'
$ws.Range("E11").Value = '  +2.44%  '

$ws.Range("E12").Value = '  +1.00%  '

$ws.Range("E13").Value = '  +3.11%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.47'
$ws.Range("E14").Value = '  +2.87%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.165.65'
$ws.Range("E15").Value = '  +4.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.386.15'
$ws.Range("E16").Value = '  +1.39%  '

$ws.Range("E17").Value = '  +1.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.671.76'
$ws.Range("E18").Value = '  +4.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.51'
$ws.Range("E19").Value = '  +4.83%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.41'
$ws.Range("E20").Value = '  +2.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '339.51'
$ws.Range("E21").Value = '  +0.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.90'
$ws.Range("E22").Value = '  +4.83%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.69'
$ws.Range("E24").Value = '  +1.46%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.66'
$ws.Range("E25").Value = '  +6.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.54'
$ws.Range("E26").Value = '  +2.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.165'
$ws.Range("E27").Value = '  +0.55%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.57'
$ws.Range("E28").Value = '  +5.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '546.40'
$ws.Range("E29").Value = '  +21.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.90'
$ws.Range("E31").Value = '  +0.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0815'
$ws.Range("E34").Value = '  +3.50%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '173.15'
$ws.Range("E35").Value = '  -1.84%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.10'
$ws.Range("E36").Value = '  +16.31%  '

$ws.Range("E37").Value = '  +2.18%  '

$ws.Range("E38").Value = '  -0.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.29'
$ws.Range("E39").Value = '  +3.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.83'
$ws.Range("E40").Value = '  +9.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '175.88'
$ws.Range("E41").Value = '  +13.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.77'
$ws.Range("E43").Value = '  +3.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.22'
$ws.Range("E44").Value = '  +7.30%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0565'
$ws.Range("E45").Value = '  +6.59%  '

$ws.Range("E46").Value = '  +0.93%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.06'
$ws.Range("E49").Value = '  +7.25%  '

$ws.Range("E50").Value = '  +5.38%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.34'
$ws.Range("E51").Value = '  -0.51%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.03'
$ws.Range("E32").Value = '  +6.39%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.82'
$ws.Range("E33").Value = '  +12.74%  '

$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0963'
$ws.Range("E47").Value = '  +0.95%  '

$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0240'
$ws.Range("E48").Value = '  +3.64%  '
